# Auto-generated edit script applying the numeric cell updates described in the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1059.8096
$ws.Range("L2").Value = 1757.75
$ws.Range("N2").Value = -1983.75
$ws.Range("J2").Value = 1757.75
$ws.Range("L87").Value = 176666.67
$ws.Range("N87").Value = -179162.67
$ws.Range("H87").Value = 182499.75
$ws.Range("J87").Value = 176666.67
$ws.Range("J90").Value = 176666.67
$ws.Range("H90").Value = 182499.75
$ws.Range("N90").Value = -542480.01
$ws.Range("L90").Value = 530000.01
$ws.Range("J116").Value = 3466.08
$ws.Range("M116").Value = 290.7777999999998
$ws.Range("H116").Value = 3382.7354
$ws.Range("L116").Value = 3466.08
$ws.Range("I116").Value = 3151.2222
$ws.Range("N116").Value = -10350.08
$ws.Range("K116").Value = 3151.2222
$ws.Range("L138").Value = 16910.2149
$ws.Range("H138").Value = 5485.8936
$ws.Range("J138").Value = 5636.7383
$ws.Range("N138").Value = -27190.2149

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("M32").Value = -4560.7583
$ws.Range("K32").Value = 4847.7583
$ws.Range("J32").Value = 28246
$ws.Range("H32").Value = 7521.843
$ws.Range("N32").Value = -28820
$ws.Range("L32").Value = 28246
$ws.Range("I32").Value = 4847.7583
$ws.Range("I45").Value = 1431.4117
$ws.Range("M45").Value = -1054.4117
$ws.Range("K45").Value = 1431.4117
$ws.Range("H45").Value = 1445.1428
$ws.Range("K61").Value = 3609.1155
$ws.Range("J61").Value = 7504.1
$ws.Range("H61").Value = 4691.0557
$ws.Range("I61").Value = 3609.1155
$ws.Range("M61").Value = -3397.1155
$ws.Range("N61").Value = -7928.1
$ws.Range("L61").Value = 7504.1
$ws.Range("H74").Value = 16174.477
$ws.Range("I74").Value = 12282.077
$ws.Range("M74").Value = -11408.077
$ws.Range("K74").Value = 12282.077
$ws.Range("K77").Value = 61410.38499999999
$ws.Range("M77").Value = -57042.38499999999
$ws.Range("H77").Value = 16174.477
$ws.Range("I77").Value = 12282.077
$ws.Range("K97").Value = 494.5
$ws.Range("M97").Value = 1.5
$ws.Range("H97").Value = 623.6818
$ws.Range("I97").Value = 494.5
$ws.Range("K102").Value = 2140.8438
$ws.Range("I102").Value = 2140.8438
$ws.Range("H102").Value = 2121.4243
$ws.Range("M102").Value = -518.8438000000001
$ws.Range("J136").Value = 7504.1
$ws.Range("L136").Value = 22512.3
$ws.Range("H136").Value = 4691.0557
$ws.Range("K136").Value = 10827.3465
$ws.Range("M136").Value = -8277.3465
$ws.Range("N136").Value = -27612.3
$ws.Range("I136").Value = 3609.1155
$ws.Range("H137").Value = 62665.832
$ws.Range("I137").Value = 49997.5
$ws.Range("M137").Value = -44897.5
$ws.Range("K137").Value = 49997.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("K105").Value = 1349.9524
$ws.Range("I105").Value = 1349.9524
$ws.Range("M105").Value = 397.0476000000001
$ws.Range("H105").Value = 1449.8462
$ws.Range("J105").Value = 1869.4
$ws.Range("L105").Value = 1869.4
$ws.Range("N105").Value = -5363.4
$ws.Range("I107").Value = 2464.7334
$ws.Range("H107").Value = 7315.3184
$ws.Range("K107").Value = 2464.7334
$ws.Range("M107").Value = -544.7334000000001
$ws.Range("H134").Value = 1563.25
$ws.Range("M134").Value = -509.1666
$ws.Range("I134").Value = 1014.7222
$ws.Range("K134").Value = 3044.1666

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("J16").Value = 2593.375
$ws.Range("L16").Value = 2593.375
$ws.Range("H16").Value = 2703.5557
$ws.Range("N16").Value = -3167.375
$ws.Range("I16").Value = 3585
$ws.Range("K16").Value = 3585
$ws.Range("M16").Value = -3298
$ws.Range("N31").Value = -883674.7
$ws.Range("H31").Value = 491921
$ws.Range("L31").Value = 883084.7
$ws.Range("K31").Value = 10488.77
$ws.Range("I31").Value = 10488.77
$ws.Range("J31").Value = 883084.7
$ws.Range("M31").Value = -10193.77
$ws.Range("M32").Value = -74
$ws.Range("K32").Value = 390
$ws.Range("J32").Value = 3000
$ws.Range("H32").Value = 1695
$ws.Range("N32").Value = -3632
$ws.Range("L32").Value = 3000
$ws.Range("I32").Value = 390
$ws.Range("N34").Value = -883488.7
$ws.Range("K34").Value = 10488.77
$ws.Range("J34").Value = 883084.7
$ws.Range("L34").Value = 883084.7
$ws.Range("H34").Value = 491921
$ws.Range("M34").Value = -10286.77
$ws.Range("I34").Value = 10488.77
$ws.Range("M58").Value = -2388.4546
$ws.Range("N58").Value = -3127.6667
$ws.Range("L58").Value = 2721.6667
$ws.Range("H58").Value = 2650.05
$ws.Range("I58").Value = 2591.4546
$ws.Range("J58").Value = 2721.6667
$ws.Range("K58").Value = 2591.4546
$ws.Range("H99").Value = 4191.8667
$ws.Range("I99").Value = 4355.643
$ws.Range("K99").Value = 4355.643
$ws.Range("J99").Value = 1899
$ws.Range("L99").Value = 1899
$ws.Range("N99").Value = -4895
$ws.Range("M99").Value = -2857.643
$ws.Range("L113").Value = 2593.375
$ws.Range("K113").Value = 3585
$ws.Range("J113").Value = 2593.375
$ws.Range("I113").Value = 3585
$ws.Range("H113").Value = 2703.5557
$ws.Range("M113").Value = -1415
$ws.Range("N113").Value = -6933.375
$ws.Range("I122").Value = 1805
$ws.Range("M122").Value = -2965
$ws.Range("K122").Value = 5415
$ws.Range("J122").Value = 2993
$ws.Range("L122").Value = 8979
$ws.Range("N122").Value = -13879
$ws.Range("H122").Value = 1947.56
$ws.Range("K126").Value = 13066.929
$ws.Range("J126").Value = 1899
$ws.Range("M126").Value = -10596.929
$ws.Range("H126").Value = 4191.8667
$ws.Range("L126").Value = 5697
$ws.Range("N126").Value = -10637
$ws.Range("I126").Value = 4355.643
$ws.Range("J136").Value = 2721.6667
$ws.Range("L136").Value = 8165.000100000001
$ws.Range("H136").Value = 2650.05
$ws.Range("K136").Value = 7774.3638
$ws.Range("M136").Value = -5224.3638
$ws.Range("N136").Value = -13265.0001
$ws.Range("I136").Value = 2591.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("M5").Value = -2950.1429
$ws.Range("K5").Value = 3062.1429
$ws.Range("L5").Value = 2181.75
$ws.Range("N5").Value = -2405.75
$ws.Range("H5").Value = 914
$ws.Range("I5").Value = 1020.7143
$ws.Range("J5").Value = 727.25
$ws.Range("I104").Value = 0
$ws.Range("H104").Value = 0
$ws.Range("M104").ClearContents()
$ws.Range("K104").Value = 0
$ws.Range("K113").Value = 11624.7276
$ws.Range("I113").Value = 3874.9092
$ws.Range("H113").Value = 22226212
$ws.Range("M113").Value = -9454.7276
$ws.Range("N127").Value = -11009.9
$ws.Range("L127").Value = 1089.9
$ws.Range("H127").Value = 363.3
$ws.Range("J127").Value = 363.3
$ws.Range("I135").Value = 1020.7143
$ws.Range("M135").Value = -6651.4287
$ws.Range("J135").Value = 727.25
$ws.Range("H135").Value = 914
$ws.Range("N135").Value = -11615.25
$ws.Range("L135").Value = 6545.25
$ws.Range("K135").Value = 9186.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("L112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L113").Value = 3089.75
$ws.Range("K113").Value = 2850.75
$ws.Range("J113").Value = 3089.75
$ws.Range("I113").Value = 2850.75
$ws.Range("H113").Value = 3030
$ws.Range("M113").Value = -680.75
$ws.Range("N113").Value = -7429.75
$ws.Range("H132").Value = 4252.8647
$ws.Range("M132").Value = -9749.75
$ws.Range("K132").Value = 12279.75
$ws.Range("I132").Value = 4093.25
$ws.Range("L141").Value = 0
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 6268
$ws.Range("H7").Value = 6955.643
$ws.Range("J7").Value = 8674.75
$ws.Range("K7").Value = 6268
$ws.Range("N7").Value = -8898.75
$ws.Range("L7").Value = 8674.75
$ws.Range("M7").Value = -6156
$ws.Range("J16").Value = 762.3333
$ws.Range("L16").Value = 762.3333
$ws.Range("H16").Value = 611.875
$ws.Range("N16").Value = -1102.3333
$ws.Range("I46").Value = 55556890
$ws.Range("M46").Value = -55556702
$ws.Range("K46").Value = 55556890
$ws.Range("H46").Value = 55556890
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("L90").Value = 0
$ws.Range("K126").Value = 18804
$ws.Range("J126").Value = 8674.75
$ws.Range("M126").Value = -16334
$ws.Range("H126").Value = 6955.643
$ws.Range("L126").Value = 26024.25
$ws.Range("N126").Value = -30964.25
$ws.Range("I126").Value = 6268

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 30025
$ws.Range("I40").Value = 30025
$ws.Range("M40").Value = -29876
$ws.Range("K40").Value = 30025
$ws.Range("I122").Value = 21742422
$ws.Range("M122").Value = -65224816
$ws.Range("K122").Value = 65227266
$ws.Range("H122").Value = 17861082
$ws.Range("H132").Value = 2693.762
$ws.Range("L132").Value = 22226.5005
$ws.Range("J132").Value = 7408.8335
$ws.Range("M132").Value = -3193.7498
$ws.Range("K132").Value = 5723.7498
$ws.Range("I132").Value = 1907.9166
$ws.Range("N132").Value = -27286.5005
$ws.Range("J136").Value = 9293.786
$ws.Range("L136").Value = 27881.358
$ws.Range("H136").Value = 6497.591
$ws.Range("K136").Value = 15578.1
$ws.Range("M136").Value = -13028.1
$ws.Range("N136").Value = -32981.358
$ws.Range("I136").Value = 5192.7
$ws.Range("L139").Value = 0
$ws.Range("H139").Value = 0
$ws.Range("N139").ClearContents()
$ws.Range("J139").Value = 0
